$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "5.20", "1.00",
# "0.0520", "62.203.27") that must stay literal text -- on a General-formatted
# cell, assigning a numeric-looking string would silently coerce it to a Number
# and drop significant trailing zeros, so force the Text format first.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.203.27"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.22"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.52"
$ws.Range("E5").Value = "  +3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.99"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.440.71"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("E11").Value = "  +2.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.895.85"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.216.47"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.443.00"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.78"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.08"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.97"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.64"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.20"
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "590.29"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0965"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.570.30"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.988"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.89"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.19"
$ws.Range("E39").Value = "  +4.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.39"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.86"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.60"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("E48").Value = "  +15.99%  "
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0521"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  -0.61%  "
